$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-27 (Timestamp, IGCC Import (MW), IGCC Export (MW))
$data = @(
    @(45821, 0, 17.029),
    @(45821.01041666666, 3.071, 1.872),
    @(45821.02083333334, 5.766, 4.823),
    @(45821.03125, 1.931, 0.21),
    @(45821.04166666666, 2.495, 0),
    @(45821.05208333334, 8.771000000000001, 0.029),
    @(45821.0625, 31.398, 0),
    @(45821.07291666666, 16.82, 0.07000000000000001),
    @(45821.08333333334, 13.029, 0.737),
    @(45821.09375, 32.645, 0),
    @(45821.10416666666, 16.655, 0.008),
    @(45821.11458333334, 13.633, 0.161),
    @(45821.125, 6.688, 2.203),
    @(45821.13541666666, 39.379, 0),
    @(45821.14583333334, 26.741, 0),
    @(45821.15625, 23.157, 0),
    @(45821.16666666666, 16.636, 0),
    @(45821.17708333334, 22.753, 0),
    @(45821.1875, 36.927, 0),
    @(45821.19791666666, 24.143, 0),
    @(45821.20833333334, 22.49, 0),
    @(45821.21875, 3.185, 2.935),
    @(45821.22916666666, 5.715, 1.527),
    @(45821.23958333334, 10.011, 0.532),
    @(45821.25, 0.098, 9.128),
    @(45821.26041666666, 1.257, 5.11)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the old trailing rows 28-42 so the used range shrinks to A1:C27
$ws.Range("A28:C42").EntireRow.Delete()
